$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking data refresh (GitHub Actions bot).
#
# The Price column (D) stores plain-text quotes (e.g. "0.9990",
# "28.729.31"). Some of the new quotes look like ordinary decimal
# numbers (e.g. "0.9990", "0.00001067") and Excel would silently
# coerce a bare Range.Value assignment into a Number, dropping the
# trailing zeros / switching to scientific notation. Pre-marking
# just those cells as Text keeps the literal string intact; the
# "thousands-dotted" quotes (e.g. "28.729.31") are never valid
# numeric literals so they do not need this treatment.
$ws.Range("D4:D15").NumberFormat = "@"
$ws.Range("D17:D22").NumberFormat = "@"
$ws.Range("D24:D28").NumberFormat = "@"
$ws.Range("D30:D31").NumberFormat = "@"
$ws.Range("D33:D40").NumberFormat = "@"
$ws.Range("D42:D46").NumberFormat = "@"
$ws.Range("D48:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.729.31"
$ws.Range("E2").Value = "  +2.42%  "
$ws.Range("D3").Value = "1.806.97"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("D4").Value = "0.9990"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "314.03"
$ws.Range("E5").Value = "  -0.53%  "
$ws.Range("D6").Value = "0.9994"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").Value = "0.5391"
$ws.Range("E7").Value = "  -1.92%  "
$ws.Range("D8").Value = "0.3791"
$ws.Range("E8").Value = "  -0.36%  "
$ws.Range("D9").Value = "0.07547"
$ws.Range("E9").Value = "  -0.58%  "
$ws.Range("D10").Value = "42.66"
$ws.Range("E10").Value = "  -1.07%  "
$ws.Range("D11").Value = "1.121"
$ws.Range("E11").Value = "  -1.32%  "
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").Value = "21.07"
$ws.Range("E13").Value = "  -0.43%  "
$ws.Range("D14").Value = "6.192"
$ws.Range("E14").Value = "  -0.62%  "
$ws.Range("D15").Value = "7.419"
$ws.Range("E15").Value = "  +3.60%  "
$ws.Range("D16").Value = "1.799.47"
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").Value = "90.78"
$ws.Range("E17").Value = "  -1.30%  "
$ws.Range("D18").Value = "0.00001067"
$ws.Range("E18").Value = "  -1.28%  "
$ws.Range("D19").Value = "0.06452"
$ws.Range("E19").Value = "  -0.69%  "
$ws.Range("D20").Value = "0.9999"
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("D21").Value = "17.31"
$ws.Range("E21").Value = "  +0.76%  "
$ws.Range("D22").Value = "5.941"
$ws.Range("E22").Value = "  -0.76%  "
$ws.Range("D23").Value = "28.680.42"
$ws.Range("E23").Value = "  +2.17%  "
$ws.Range("D24").Value = "11.22"
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("D25").Value = "2.113"
$ws.Range("E25").Value = "  +1.10%  "
$ws.Range("D26").Value = "160.86"
$ws.Range("E26").Value = "  +3.24%  "
$ws.Range("D27").Value = "20.61"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").Value = "2.386"
$ws.Range("E28").Value = "  -0.11%  "
$ws.Range("D29").Value = "2.004.37"
$ws.Range("E29").Value = "  -0.24%  "
$ws.Range("D30").Value = "123.82"
$ws.Range("E30").Value = "  +0.67%  "
$ws.Range("D31").Value = "1.115"
$ws.Range("E31").Value = "  -3.28%  "
$ws.Range("E32").Value = "  +0.33%  "
$ws.Range("D33").Value = "5.695"
$ws.Range("E33").Value = "  -0.97%  "
$ws.Range("D34").Value = "3.699"
$ws.Range("E34").Value = "  +2.73%  "
$ws.Range("D35").Value = "0.2266"
$ws.Range("E35").Value = "  +6.56%  "
$ws.Range("D36").Value = "0.06509"
$ws.Range("E36").Value = "  +7.65%  "
$ws.Range("D37").Value = "8.949"
$ws.Range("E37").Value = "  +3.26%  "
$ws.Range("D38").Value = "0.02319"
$ws.Range("E38").Value = "  +0.60%  "
$ws.Range("D39").Value = "5.063"
$ws.Range("E39").Value = "  +0.79%  "
$ws.Range("D40").Value = "11.38"
$ws.Range("E40").Value = "  -1.37%  "
$ws.Range("E41").Value = "  -0.83%  "
$ws.Range("D42").Value = "1.208"
$ws.Range("E42").Value = "  +4.87%  "
$ws.Range("D43").Value = "0.9994"
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").Value = "1.393"
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("D45").Value = "13.35"
$ws.Range("E45").Value = "  -0.81%  "
$ws.Range("D46").Value = "0.5898"
$ws.Range("E46").Value = "  -0.48%  "
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("D48").Value = "126.38"
$ws.Range("E48").Value = "  +3.66%  "
$ws.Range("D49").Value = "1.966"
$ws.Range("E49").Value = "  +1.95%  "
$ws.Range("D50").Value = "1.159"
$ws.Range("E50").Value = "  +2.03%  "
$ws.Range("D51").Value = "0.06892"
$ws.Range("E51").Value = "  +1.49%  "
